$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '96.836.75'
$ws.Range('E2').Value = '  +1.36%  '

$ws.Range('D3').Value = '3.582.40'
$ws.Range('E3').Value = '  +0.67%  '

$ws.Range('E4').Value = '  -0.02%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '240.24'
$c.NumberFormat = 'General'
$ws.Range('E5').Value = '  +1.44%  '

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '1.79'
$c.NumberFormat = 'General'
$ws.Range('E6').Value = '  +17.08%  '

$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '649.47'
$c.NumberFormat = 'General'
$ws.Range('E7').Value = '  -0.47%  '

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.419'
$c.NumberFormat = 'General'
$ws.Range('E8').Value = '  +4.59%  '

$ws.Range('B9').Value = 'USDC'
$ws.Range('C9').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.NumberFormat = 'General'
$ws.Range('E9').Value = '  -0.03%  '

$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '1.06'
$c.NumberFormat = 'General'
$ws.Range('E10').Value = '  +2.16%  '

$ws.Range('D11').Value = '3.581.36'
$ws.Range('E11').Value = '  +0.70%  '

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '44.19'
$c.NumberFormat = 'General'
$ws.Range('E12').Value = '  +3.17%  '

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.203'
$c.NumberFormat = 'General'
$ws.Range('E13').Value = '  +0.40%  '

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '6.45'
$c.NumberFormat = 'General'
$ws.Range('E14').Value = '  +1.01%  '

$ws.Range('D15').Value = '4.248.49'
$ws.Range('E15').Value = '  +0.63%  '

$ws.Range('D16').Value = '96.821.84'
$ws.Range('E16').Value = '  +1.51%  '

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.0000257'
$c.NumberFormat = 'General'
$ws.Range('E17').Value = '  +1.17%  '

$ws.Range('D18').Value = '3.548.73'

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '8.74'
$c.NumberFormat = 'General'
$ws.Range('E19').Value = '  +2.37%  '

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '12.71'
$c.NumberFormat = 'General'
$ws.Range('E20').Value = '  +0.47%  '

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '18.07'
$c.NumberFormat = 'General'
$ws.Range('E21').Value = '  +1.51%  '

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '0.533'
$c.NumberFormat = 'General'
$ws.Range('E22').Value = '  +9.05%  '

$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '513.34'
$c.NumberFormat = 'General'
$ws.Range('E23').Value = '  +0.84%  '

$ws.Range('B24').Value = 'SuiNetwork'
$ws.Range('C24').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '3.43'
$c.NumberFormat = 'General'
$ws.Range('E24').Value = '  -0.21%  '

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '0.0000203'
$c.NumberFormat = 'General'
$ws.Range('E25').Value = '  +1.46%  '

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '6.86'
$c.NumberFormat = 'General'
$ws.Range('E26').Value = '  -0.32%  '

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '102.42'
$c.NumberFormat = 'General'
$ws.Range('E27').Value = '  +6.92%  '

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '13.31'
$c.NumberFormat = 'General'
$ws.Range('E28').Value = '  +4.86%  '

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '0.173'
$c.NumberFormat = 'General'
$ws.Range('E29').Value = '  +19.70%  '

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '2.99'
$c.NumberFormat = 'General'
$ws.Range('E30').Value = '  -0.94%  '

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '11.88'
$c.NumberFormat = 'General'
$ws.Range('E31').Value = '  +2.97%  '

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.NumberFormat = 'General'
$ws.Range('E32').Value = '  +0.05%  '

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.186'
$c.NumberFormat = 'General'
$ws.Range('E33').Value = '  +4.90%  '

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.NumberFormat = 'General'
$ws.Range('E34').Value = '  +0.10%  '

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '31.84'
$c.NumberFormat = 'General'
$ws.Range('E35').Value = '  +1.33%  '

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.73'
$c.NumberFormat = 'General'
$ws.Range('E36').Value = '  +6.82%  '

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.579'
$c.NumberFormat = 'General'
$ws.Range('E37').Value = '  +2.42%  '

$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '8.78'
$c.NumberFormat = 'General'
$ws.Range('E38').Value = '  +4.11%  '

$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '611.39'
$c.NumberFormat = 'General'
$ws.Range('E39').Value = '  +3.75%  '

$ws.Range('E40').Value = '  +2.44%  '

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '1.91'
$c.NumberFormat = 'General'
$ws.Range('E41').Value = '  +3.49%  '

$ws.Range('E42').Value = '  -0.02%  '

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.926'
$c.NumberFormat = 'General'
$ws.Range('E43').Value = '  +2.98%  '

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '6.15'
$c.NumberFormat = 'General'
$ws.Range('E44').Value = '  +6.56%  '

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.0452'
$c.NumberFormat = 'General'
$ws.Range('E45').Value = '  +9.00%  '

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.426'
$c.NumberFormat = 'General'
$ws.Range('E46').Value = '  +27.79%  '

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '2.28'
$c.NumberFormat = 'General'
$ws.Range('E47').Value = '  -0.12%  '

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '23.61'
$c.NumberFormat = 'General'
$ws.Range('E48').Value = '  +0.73%  '

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '8.64'
$c.NumberFormat = 'General'
$ws.Range('E49').Value = '  +4.55%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '33.40'
$c.NumberFormat = 'General'
$ws.Range('E50').Value = '  -1.52%  '

$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '3.24'
$c.NumberFormat = 'General'
$ws.Range('E51').Value = '  +4.90%  '
